# "Add files via upload" - the player roster backing this sheet was
# refreshed: one new player (Mike Conley / Minnesota Timberwolves) was
# added to the source data, and the full table was re-uploaded, changing
# the row order for the existing players too.
#
# Rebuild the data block (A2:C18) to match the newly uploaded table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Shai Gilgeous-Alexander", "PG",          "Oklahoma City Thunder"),
    @("CJ McCollum",             "PG,SG",       "New Orleans Pelicans"),
    @("Jordan Poole",            "PG,SG",       "Washington Wizards"),
    @("Kyrie Irving",            "PG,SG",       "Dallas Mavericks"),
    @("RJ Barrett",              "SF,PF",       "Toronto Raptors"),
    @("Tobias Harris",           "SF,PF",       "Detroit Pistons"),
    @("Brandon Boston Jr.",      "SG,SF,PF",    "New Orleans Pelicans"),
    @("Christian Braun",         "SG,SF",       "Denver Nuggets"),
    @("Jalen Williams",          "SG,SF,PF,C",  "Oklahoma City Thunder"),
    @("Jimmy Butler",            "SF,PF",       "Miami Heat"),
    @("Mike Conley",             "PG",          "Minnesota Timberwolves"),
    @("John Collins",            "PF,C",        "Utah Jazz"),
    @("Dennis Schröder",         "PG",          "Brooklyn Nets"),
    @("Keyonte George",          "PG,SG",       "Utah Jazz"),
    @("Zach LaVine",             "SG,SF",       "Chicago Bulls"),
    @("Joel Embiid",             "C",           "Philadelphia 76ers"),
    @("Lauri Markkanen",         "SF,PF",       "Utah Jazz")
)

# Clear the old table body first (header row 1 is untouched).
$ws.Range("A2:C17").ClearContents()

$row = 2
foreach ($player in $data) {
    $ws.Cells.Item($row, 1).Value = $player[0]
    $ws.Cells.Item($row, 2).Value = $player[1]
    $ws.Cells.Item($row, 3).Value = $player[2]
    $row++
}
